$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShape($inlineShape, $newName) {
    # InlineShape has no settable .Name in the Word object model; the
    # <wp:docPr name="..."> is only reachable through the floating Shape
    # object, so round-trip through ConvertToShape/ConvertToInlineShape.
    $floating = $inlineShape.ConvertToShape()
    $floating.Name = $newName
    $floating.ConvertToInlineShape() | Out-Null
}

# Footer (default / odd pages) - Pearson logo, wp:docPr id="1"
$footerDefault = $sec.Footers.Item(1)
$shape1 = $footerDefault.Range.InlineShapes.Item(1)
Rename-InlineShape $shape1 "image1.png"

# Footer (first page) - Pearson logo, wp:docPr id="2"
$footerFirst = $sec.Footers.Item(2)
$shape2 = $footerFirst.Range.InlineShapes.Item(1)
Rename-InlineShape $shape2 "image1.png"

# Header (first page) - BTEC logo, wp:docPr id="3"
$headerFirst = $sec.Headers.Item(2)
$shape3 = $headerFirst.Range.InlineShapes.Item(1)
Rename-InlineShape $shape3 "image2.jpg"
